$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 16 (pushes the old "empty spacer" row 15 pattern
# down to 16, and the totals block 16-18 down to 17-19)
$ws.Rows("16").Insert()

# Fill in the newly available row 15 with the latest timeline entry
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "26/5/2024"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = "Finished search, pagination, recommended products, my products, favorite products"

# Extend the total-hours formula to include the new row
$ws.Range("D17").Formula = "=SUM(C4:C15)"

# The engine auto-expands the row height of row 17 because of the large
# font used by the totals styling; restore it to the default (no explicit
# row height) like the rest of the untouched rows.
$ws.Rows("17").AutoFit()

# Update the saved selection to mirror what Excel would have left it at
# after editing/extending that area.
$ws.Range("D17:D19").Select()
